$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '72.855.14'
$ws.Range("E2").Value = '  +1.18%  '
Set-TextValue $ws.Range("D3") '3.982.04'
$ws.Range("E3").Value = '  -0.96%  '
Set-TextValue $ws.Range("D4") '0.998'
$ws.Range("E4").Value = '  -0.20%  '
Set-TextValue $ws.Range("D5") '617.31'
$ws.Range("E5").Value = '  +15.38%  '
Set-TextValue $ws.Range("D6") '166.24'
$ws.Range("E6").Value = '  +10.59%  '
Set-TextValue $ws.Range("D7") '0.686'
$ws.Range("E7").Value = '  -1.25%  '
$ws.Range("E8").Value = '  -0.06%  '
Set-TextValue $ws.Range("D9") '0.759'
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("E10").Value = '  -2.04%  '
Set-TextValue $ws.Range("D11") '57.66'
$ws.Range("E11").Value = '  +6.25%  '
Set-TextValue $ws.Range("D12") '0.0000316'
$ws.Range("E12").Value = '  -3.00%  '
Set-TextValue $ws.Range("D13") '11.20'
$ws.Range("E13").Value = '  +3.57%  '
Set-TextValue $ws.Range("D14") '4.605.00'
$ws.Range("E14").Value = '  -1.25%  '
Set-TextValue $ws.Range("D15") '3.980.98'
$ws.Range("E15").Value = '  -0.91%  '
Set-TextValue $ws.Range("D16") '1.26'
$ws.Range("E16").Value = '  +5.95%  '
Set-TextValue $ws.Range("D17") '14.36'
$ws.Range("E17").Value = '  +1.22%  '
Set-TextValue $ws.Range("D18") '20.79'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("E19").Value = '  +0.22%  '
Set-TextValue $ws.Range("D20") '72.705.43'
$ws.Range("E20").Value = '  +1.00%  '
Set-TextValue $ws.Range("D21") '442.52'
$ws.Range("E21").Value = '  +2.28%  '
Set-TextValue $ws.Range("D22") '4.97'
$ws.Range("E22").Value = '  +17.54%  '
Set-TextValue $ws.Range("D23") '96.55'
$ws.Range("E23").Value = '  -2.13%  '
Set-TextValue $ws.Range("D24") '3.40'
$ws.Range("E24").Value = '  -4.90%  '
Set-TextValue $ws.Range("D25") '14.60'
$ws.Range("E25").Value = '  -0.94%  '
Set-TextValue $ws.Range("D26") '4.28'
$ws.Range("E26").Value = '  -1.01%  '
Set-TextValue $ws.Range("D27") '11.39'
$ws.Range("E27").Value = '  -0.91%  '
Set-TextValue $ws.Range("D28") '10.55'
$ws.Range("E28").Value = '  -2.73%  '
Set-TextValue $ws.Range("D29") '5.93'
$ws.Range("E29").Value = '  +0.27%  '
Set-TextValue $ws.Range("D30") '36.23'
$ws.Range("E30").Value = '  -2.11%  '
Set-TextValue $ws.Range("D31") '7.84'
$ws.Range("E31").Value = '  -4.65%  '
Set-TextValue $ws.Range("D32") '13.92'
$ws.Range("E32").Value = '  +2.52%  '
Set-TextValue $ws.Range("D33") '0.131'
$ws.Range("E33").Value = '  -3.94%  '
Set-TextValue $ws.Range("D34") '48.60'
$ws.Range("E34").Value = '  -3.26%  '
Set-TextValue $ws.Range("D35") '71.96'
$ws.Range("E35").Value = '  +10.04%  '
Set-TextValue $ws.Range("D36") '640.74'
$ws.Range("E36").Value = '  -5.75%  '
Set-TextValue $ws.Range("D37") '0.0₃0896'
$ws.Range("E37").Value = '  +7.70%  '
Set-TextValue $ws.Range("D38") '0.437'
$ws.Range("E38").Value = '  -4.23%  '
Set-TextValue $ws.Range("D39") '3.50'
$ws.Range("E39").Value = '  +3.40%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D40") '0.149'
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("B41").Value = 'Dai'
$ws.Range("C41").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D41") '0.999'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D42") '3.35'
$ws.Range("E42").Value = '  -1.34%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D44") '10.85'
$ws.Range("E44").Value = '  -1.23%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D45") '0.0488'
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D47") '2.66'
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D48") '3.44'
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D49") '2.922.90'
$ws.Range("E49").Value = '  +2.60%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D50") '3.10'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D51") '2.83'
$ws.Range("E51").Value = '  +29.92%  '
